# Insert a new data row at row 108 (pushing the existing rows 108-127 down to
# 109-128), and populate it with the new "Ají" / "Cristal" price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(108).Insert()

$ws.Cells.Item(108, 1).Value2 = 11
$ws.Cells.Item(108, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(108, 3).Value2 = "Bíobío"
$ws.Cells.Item(108, 4).Value2 = 44798
$ws.Cells.Item(108, 5).Value2 = 8
$ws.Cells.Item(108, 6).Value2 = 100112021
$ws.Cells.Item(108, 7).Value2 = "Ají"
$ws.Cells.Item(108, 8).Value2 = "Cristal"
$ws.Cells.Item(108, 9).Value2 = "Primera"
$ws.Cells.Item(108, 10).Value2 = 60
$ws.Cells.Item(108, 11).Value2 = 54000
$ws.Cells.Item(108, 12).Value2 = 56000
$ws.Cells.Item(108, 13).Value2 = 55000
$ws.Cells.Item(108, 14).Value2 = "`$/caja 25 kilos"
$ws.Cells.Item(108, 15).Value2 = "Provincia de Limarí"
$ws.Cells.Item(108, 16).Value2 = 2200
$ws.Cells.Item(108, 17).Value2 = 25
$ws.Cells.Item(108, 18).Value2 = "Hortaliza"

# Keep the date cell formatted the same way as the other date cells in column D.
$ws.Cells.Item(108, 4).NumberFormat = $ws.Cells.Item(109, 4).NumberFormat
